$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 'A partnership that has its principal office and place of business in any place not subject to the jurisdiction of the United States.'
$ws.Range("E3").Value = 0

$ws.Range("C4").Value = 'A general partner of an investment adviser who resides in any place not subject to the jurisdiction of the United States.'
$ws.Range("E4").Value = 0

$ws.Range("C5").Value = 'Each party named in the process, pleadings, or papers served on the Commission.'
$ws.Range("E5").Value = 0

$ws.Range("C6").Value = 'An individual, corporation, partnership, or unincorporated organization that resides or has its principal office and place of business in any place not subject to the jurisdiction of the United States.'
$ws.Range("E6").Value = 0.4195011337868481

$ws.Range("C7").Value = 'Documents served on a non-resident investment adviser, general partner, or managing agent.'
$ws.Range("E7").Value = 0

$ws.Range("C8").Value = 'A location not subject to the jurisdiction of the United States.'
$ws.Range("E8").Value = 0

$ws.Range("C9").Value = 'An individual or entity involved in serving process, pleadings, or papers.'
$ws.Range("E9").Value = 0

$ws.Range("C10").Value = 'An association that has its principal office and place of business in any place not subject to the jurisdiction of the United States.'
$ws.Range("E10").Value = 0

$ws.Range("C11").Value = 'A managing agent of an investment adviser who resides in any place not subject to the jurisdiction of the United States.'
$ws.Range("E11").Value = 0

$ws.Range("C12").Value = 'The records maintained by the Commission where an additional copy of the process, pleadings, or papers is kept.'
$ws.Range("E12").Value = 0

$ws.Range("C13").Value = 'The Secretary of the Commission (Secretary) will promptly forward a copy to each named party by registered or certified mail at that party''s last address filed with the Commission.'
$ws.Range("E13").Value = 0

$ws.Range("C14").Value = 'The section where a person may serve a non-resident investment adviser, general partner, or managing agent by furnishing the Commission with copies of the documents.'
$ws.Range("E14").Value = 0

$ws.Range("C15").Value = 'The main office of a corporation, partnership, or unincorporated organization.'
$ws.Range("E15").Value = 0

$ws.Range("C16").Value = 'The section that defines ''Principal office and place of business''.'
$ws.Range("E16").Value = 0

$ws.Range("C17").Value = 'Certification by the Secretary that the Commission was served with documents and forwarded them to a named party.'
$ws.Range("E17").Value = 0

$ws.Range("C18").Value = 'A person who directs or manages the affairs of any unincorporated organization or association other than a partnership.'
$ws.Range("E18").Value = 0

$ws.Range("C19").Value = 'A corporation that is incorporated in or that has its principal office and place of business in any place not subject to the jurisdiction of the United States.'
$ws.Range("E19").Value = 0

$ws.Range("C20").Value = 'An organization that is not incorporated and has its principal office and place of business in any place not subject to the jurisdiction of the United States.'
$ws.Range("E20").Value = 0

$ws.Range("C21").Value = 'A method used by the Secretary to forward documents to named parties.'
$ws.Range("E21").Value = 0

$ws.Range("C22").Value = 'An investment adviser who resides in any place not subject to the jurisdiction of the United States.'
$ws.Range("E22").Value = 0

$ws.Range("C23").Value = 'The act of the Secretary certifying that the Commission was served with documents and forwarded them to a named party.'
$ws.Range("E23").Value = 0

$ws.Range("C24").Value = 'The section where the Secretary of the Commission forwards a copy to each named party by registered or certified mail.'
$ws.Range("E24").Value = 0

$ws.Range("C25").Value = 'The location where a corporation, partnership, or unincorporated organization conducts its business.'
$ws.Range("E25").Value = 0

$ws.Range("C26").Value = 'The legal authority or control of the United States over a place.'
$ws.Range("E26").Value = 0

$ws.Range("C27").Value = 'Activities directed or managed by a managing agent of any unincorporated organization or association other than a partnership.'
$ws.Range("E27").Value = 0

$ws.Range("C28").Value = 'Agents appointed to receive service of process, pleadings, or papers on behalf of a non-resident investment adviser, general partner, or managing agent.'
$ws.Range("E28").Value = 0

$ws.Range("C29").Value = 'An individual who resides in any place not subject to the jurisdiction of the United States.'
$ws.Range("E29").Value = 0

$ws.Range("C30").Value = 'Legal documents served on a non-resident investment adviser, general partner, or managing agent.'
$ws.Range("E30").Value = 0

$ws.Range("C31").Value = 'Documents served on a non-resident investment adviser, general partner, or managing agent.'
$ws.Range("E31").Value = 0

$ws.Range("C32").Value = 'The Secretary of the Commission who forwards documents to named parties.'
$ws.Range("E32").Value = 0

$ws.Range("C33").Value = 'An adviser who provides investment advice and may be served with process, pleadings, or papers.'
$ws.Range("E33").Value = 0

$ws.Range("C34").Value = 'The last address filed with the Commission by a named party.'
$ws.Range("E34").Value = 0

$ws.Range("C35").Value = 'A method used by the Secretary to forward documents to named parties.'
$ws.Range("E35").Value = 0

$ws.Range("C36").Value = 'The entity to which a person may furnish copies of process, pleadings, or papers for service on non-resident investment advisers, general partners, or managing agents.'
$ws.Range("E36").Value = 0

$ws.Range("C37").Value = 'An order disposing of the matter will be issued following the expiration of the period of time.'
$ws.Range("E37").Value = 0.4357798165137615

$ws.Range("C38").Value = 'Reasons for requesting a hearing must be stated by the interested person.'
$ws.Range("E38").Value = 0.1857585139318886

$ws.Range("C39").Value = 'Refers to the Act under which the Commission operates and makes orders.'
$ws.Range("E39").Value = 0.25

$ws.Range("C41").Value = 'A reason for which the Commission may order a hearing.'
$ws.Range("E41").Value = 0.2048192771084337

$ws.Range("C42").Value = 'An interested person may submit facts and request a hearing, stating reasons and the nature of their interest.'
$ws.Range("E42").Value = 0.3504672897196262

$ws.Range("C43").Value = 'The earliest date upon which an order disposing of the matter may be entered, as indicated in the notice.'
$ws.Range("E43").Value = 0.2888888888888889

$ws.Range("C44").Value = 'The publication where notice of the initiation of the proceeding will be published.'
$ws.Range("E44").Value = 0.4470588235294117

$ws.Range("C45").Value = 'The initiation of the proceeding is marked by the filing of an application or upon the Commission''s own motion.'
$ws.Range("E45").Value = 0.3423423423423423

$ws.Range("C46").Value = 'Notice of the initiation of the proceeding will be published in the Federal Register.'
$ws.Range("E46").Value = 0.4722222222222222

$ws.Range("C47").Value = 'A hearing may be ordered if necessary or appropriate in the public interest or for the protection of investors.'
$ws.Range("E47").Value = 0.6144578313253012

$ws.Range("C48").Value = 'A factor considered by the Commission to determine the necessity of a hearing.'
$ws.Range("E48").Value = 0.3053435114503816

$ws.Range("C50").Value = 'Facts bearing upon the desirability of a hearing on the matter may be submitted by an interested person.'
$ws.Range("E50").Value = 0.2742857142857142

$ws.Range("C51").Value = 'The period after which an order disposing of the matter will be issued unless a hearing is ordered.'
$ws.Range("E51").Value = 0.3006535947712419

$ws.Range("C52").Value = 'An application means any application for an order of the Commission under the Act other than an application for registration as an investment adviser.'
$ws.Range("E52").Value = 0.8324022346368716

$ws.Range("C53").Value = 'An application for registration as an investment adviser is excluded from the definition of ''application''.'
$ws.Range("E53").Value = 0.5377358490566038

$ws.Range("C54").Value = 'The specified duration within which actions must be taken or submissions made.'
$ws.Range("E54").Value = 0.2514285714285714

$ws.Range("C55").Value = 'The body that initiates proceedings, issues orders, and may order hearings.'
$ws.Range("E55").Value = 0.2710843373493976

$ws.Range("C56").Value = 'Consideration for whether a hearing should be held, based on submitted facts.'
$ws.Range("E56").Value = 0.2457142857142857

$ws.Range("C57").Value = 'The subject of the proceeding initiated by the application or Commission''s motion.'
$ws.Range("E57").Value = 0.2810457516339869

$ws.Range("C58").ClearContents()
$ws.Range("E58").Value = 0

$ws.Range("C60").ClearContents()
$ws.Range("E60").Value = 0

$ws.Range("C64").ClearContents()
$ws.Range("E64").Value = 0

$ws.Range("C65").ClearContents()
$ws.Range("E65").Value = 0

$ws.Range("C68").ClearContents()
$ws.Range("E68").Value = 0

$ws.Range("C69").Value = 'Assets under management, as defined under Section 203A(a)(3) of the Act and reported on its annual updating amendment to Form ADV, of less than $25 million, or such higher amount as the Commission may by rule deem appropriate.'
$ws.Range("E69").Value = 0.8528301886792453

$ws.Range("C71").ClearContents()
$ws.Range("E71").Value = 0

$ws.Range("C72").ClearContents()
$ws.Range("E72").Value = 0

$ws.Range("C73").Value = 'The total assets as shown on the balance sheet of the investment adviser or other person, or the balance sheet of the investment adviser or such other person with its subsidiaries consolidated, whichever is larger.'
$ws.Range("E73").Value = 0.9363636363636364

$ws.Range("C75").ClearContents()
$ws.Range("E75").Value = 0

$ws.Range("C76").ClearContents()
$ws.Range("E76").Value = 0

$ws.Range("C77").Value = 'An investment adviser that has assets under management of less than $25 million, or such higher amount as the Commission may by rule deem appropriate, did not have total assets of $5 million or more on the last day of the most recent fiscal year, and does not control, is not controlled by, and is not under common control with another investment adviser that has assets under management of $25 million or more, or any person (other than a natural person) that had total assets of $5 million or more on the last day of the most recent fiscal year.'
$ws.Range("E77").Value = 0.8720292504570384

$ws.Range("C78").ClearContents()
$ws.Range("E78").Value = 0

$ws.Range("C82").ClearContents()
$ws.Range("E82").Value = 0

Write-Output "Edit applied successfully"
